# Correção nos dados e início da análise PNAD 2009
#
# The original sheet had two "section header" rows (row 5: "situação do
# domicílio" and row 8: "grandes regiões e unidades da federação") that
# carried no data - they were left over from a pandas multi-index export
# and only occupied column A, with the rest of the row blank. This edit
# removes those two empty/placeholder rows (all subsequent rows shift up
# by two), and also fixes the second header row, where the
# "unnamed: 1_level_1" / "unnamed: 5_level_1" placeholder labels coming
# from the unnamed pandas column levels are corrected to "total" to match
# the other "total" column header already present in the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two label-only rows that have no associated data.
# Row 8 ("grandes regiões e unidades da federação") is removed first so
# that row 5 ("situação do domicílio") keeps its row number "5" until it,
# too, is deleted right after.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()

# Fix the mis-labelled "unnamed" headers on row 2.
$ws.Range("B2").Value2 = "total"
$ws.Range("F2").Value2 = "total"
